$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Attribute Group" (J) and "Attribute Value" (K) columns are no longer
# needed, so remove them entirely. This shifts Product Category / Vendor /
# Sub Category Name left by two columns.
$ws.Range("J1:K1").EntireColumn.Delete()

# Update the active selection to match the post-edit state.
$ws.Range("H14").Select()
